$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4035.9443
$ws.Range("I11").Value = 4035.9443
$ws.Range("K11").Value = 4035.9443
$ws.Range("M11").Value = -3895.9443
$ws.Range("H18").Value = 914.1875
$ws.Range("I18").Value = 809.0714
$ws.Range("J18").Value = 1650
$ws.Range("K18").Value = 809.0714
$ws.Range("L18").Value = 1650
$ws.Range("M18").Value = -525.0714
$ws.Range("N18").Value = -2218
$ws.Range("H51").Value = 26333.334
$ws.Range("I51").Value = 12700
$ws.Range("J51").Value = 36071.43
$ws.Range("K51").Value = 12700
$ws.Range("L51").Value = 36071.43
$ws.Range("M51").Value = -12216
$ws.Range("N51").Value = -37039.43
$ws.Range("H55").Value = 320.875
$ws.Range("I55").Value = 348.41666
$ws.Range("J55").Value = 238.25
$ws.Range("K55").Value = 348.41666
$ws.Range("L55").Value = 238.25
$ws.Range("M55").Value = -134.41666
$ws.Range("N55").Value = -666.25
$ws.Range("H70").Value = 4183.85
$ws.Range("I70").Value = 3460.2
$ws.Range("J70").Value = 4425.067
$ws.Range("K70").Value = 10380.6
$ws.Range("L70").Value = 13275.201
$ws.Range("M70").Value = -10110.6
$ws.Range("N70").Value = -13815.201
$ws.Range("H73").Value = 4183.85
$ws.Range("I73").Value = 3460.2
$ws.Range("J73").Value = 4425.067
$ws.Range("K73").Value = 10380.6
$ws.Range("L73").Value = 13275.201
$ws.Range("M73").Value = -9444.599999999999
$ws.Range("N73").Value = -15147.201
$ws.Range("H86").Value = 2376.4
$ws.Range("I86").Value = 2511.6667
$ws.Range("J86").Value = 2173.5
$ws.Range("K86").Value = 2511.6667
$ws.Range("L86").Value = 2173.5
$ws.Range("M86").Value = -1388.6667
$ws.Range("N86").Value = -4419.5
$ws.Range("H89").Value = 2376.4
$ws.Range("I89").Value = 2511.6667
$ws.Range("J89").Value = 2173.5
$ws.Range("K89").Value = 12558.3335
$ws.Range("L89").Value = 10867.5
$ws.Range("M89").Value = -6942.333500000001
$ws.Range("N89").Value = -22099.5
$ws.Range("H100").Value = 3014.8572
$ws.Range("I100").Value = 1601.6666
$ws.Range("J100").Value = 4074.75
$ws.Range("K100").Value = 1601.6666
$ws.Range("L100").Value = 4074.75
$ws.Range("M100").Value = -1060.6666
$ws.Range("N100").Value = -5156.75
$ws.Range("H101").Value = 1636.375
$ws.Range("I101").Value = 1656
$ws.Range("J101").Value = 1499
$ws.Range("K101").Value = 4968
$ws.Range("L101").Value = 4497
$ws.Range("M101").Value = -3346
$ws.Range("N101").Value = -7741
$ws.Range("H112").Value = 3302.75
$ws.Range("J112").Value = 3302.75
$ws.Range("L112").Value = 9908.25
$ws.Range("N112").Value = -12124.25
$ws.Range("H132").Value = 1702.9231
$ws.Range("I132").Value = 1469.5
$ws.Range("K132").Value = 4408.5
$ws.Range("M132").Value = -1878.5
$ws.Range("H138").Value = 3128.42
$ws.Range("I138").Value = 2113.3
$ws.Range("J138").Value = 3382.2
$ws.Range("K138").Value = 6339.900000000001
$ws.Range("L138").Value = 10146.6
$ws.Range("M138").Value = -1199.900000000001
$ws.Range("N138").Value = -20426.6
$ws.Range("H141").Value = 830.44446
$ws.Range("I141").Value = 830.44446
$ws.Range("K141").Value = 2491.33338
$ws.Range("M141").Value = 2688.66662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8502.733
$ws.Range("I45").Value = 9349.308000000001
$ws.Range("K45").Value = 9349.308000000001
$ws.Range("M45").Value = -8972.308000000001
$ws.Range("H50").Value = 2074.2222
$ws.Range("I50").Value = 5300
$ws.Range("J50").Value = 1152.5714
$ws.Range("K50").Value = 5300
$ws.Range("L50").Value = 1152.5714
$ws.Range("M50").Value = -4586
$ws.Range("N50").Value = -2580.5714
$ws.Range("H61").Value = 5899.154
$ws.Range("I61").Value = 3404
$ws.Range("K61").Value = 3404
$ws.Range("M61").Value = -3192
$ws.Range("H74").Value = 388614.38
$ws.Range("I74").Value = 437824.97
$ws.Range("K74").Value = 437824.97
$ws.Range("M74").Value = -436950.97
$ws.Range("H77").Value = 388614.38
$ws.Range("I77").Value = 437824.97
$ws.Range("K77").Value = 2189124.85
$ws.Range("M77").Value = -2184756.85
$ws.Range("H122").Value = 976.08
$ws.Range("I122").Value = 1078.2222
$ws.Range("K122").Value = 3234.6666
$ws.Range("M122").Value = -784.6665999999996
$ws.Range("H136").Value = 5899.154
$ws.Range("I136").Value = 3404
$ws.Range("K136").Value = 10212
$ws.Range("M136").Value = -7662

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 258000
$ws.Range("I26").Value = 258000
$ws.Range("K26").Value = 258000
$ws.Range("M26").Value = -257708
$ws.Range("H99").Value = 3574.375
$ws.Range("I99").Value = 3087.25
$ws.Range("K99").Value = 3087.25
$ws.Range("M99").Value = -1589.25
$ws.Range("H105").Value = 27035418
$ws.Range("I105").Value = 35724732
$ws.Range("K105").Value = 35724732
$ws.Range("M105").Value = -35722985
$ws.Range("H134").Value = 6089.1665
$ws.Range("I134").Value = 2587.9092
$ws.Range("K134").Value = 7763.7276
$ws.Range("M134").Value = -5228.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 806.1818
$ws.Range("I22").Value = 637.2857
$ws.Range("J22").Value = 1101.75
$ws.Range("K22").Value = 637.2857
$ws.Range("L22").Value = 1101.75
$ws.Range("M22").Value = -287.2857
$ws.Range("N22").Value = -1801.75
$ws.Range("H31").Value = 29414782
$ws.Range("I31").Value = 41667896
$ws.Range("J31").Value = 7311.9
$ws.Range("K31").Value = 41667896
$ws.Range("L31").Value = 7311.9
$ws.Range("M31").Value = -41667601
$ws.Range("N31").Value = -7901.9
$ws.Range("H34").Value = 29414782
$ws.Range("I34").Value = 41667896
$ws.Range("J34").Value = 7311.9
$ws.Range("K34").Value = 41667896
$ws.Range("L34").Value = 7311.9
$ws.Range("M34").Value = -41667694
$ws.Range("N34").Value = -7715.9
$ws.Range("H37").Value = 22666.666
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 22666.666
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 22666.666
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -22880.666
$ws.Range("H107").Value = 444.25
$ws.Range("I107").Value = 442.33334
$ws.Range("K107").Value = 442.33334
$ws.Range("M107").Value = 1477.66666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1860784
$ws.Range("I4").Value = 2962012
$ws.Range("J4").Value = 81877
$ws.Range("K4").Value = 8886036
$ws.Range("L4").Value = 245631
$ws.Range("M4").Value = -8885924
$ws.Range("N4").Value = -245855
$ws.Range("H34").Value = 2058.3
$ws.Range("I34").Value = 58.714287
$ws.Range("J34").Value = 3135
$ws.Range("K34").Value = 176.142861
$ws.Range("L34").Value = 9405
$ws.Range("M34").Value = -92.14286099999998
$ws.Range("N34").Value = -9573
$ws.Range("H106").Value = 6666.3335
$ws.Range("I106").Value = 5000
$ws.Range("J106").Value = 7499.5
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 22498.5
$ws.Range("M106").Value = -14054
$ws.Range("N106").Value = -24390.5
$ws.Range("H131").Value = 9012489
$ws.Range("J131").Value = 3982.5483
$ws.Range("L131").Value = 11947.6449
$ws.Range("N131").Value = -22027.6449

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 19776
$ws.Range("I113").Value = 19776
$ws.Range("K113").Value = 19776
$ws.Range("M113").Value = -17606
$ws.Range("H132").Value = 4214.826
$ws.Range("I132").Value = 2872.8125
$ws.Range("J132").Value = 7282.2856
$ws.Range("K132").Value = 8618.4375
$ws.Range("L132").Value = 21846.8568
$ws.Range("M132").Value = -6088.4375
$ws.Range("N132").Value = -26906.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6861.9
$ws.Range("J46").Value = 7775.7646
$ws.Range("L46").Value = 7775.7646
$ws.Range("N46").Value = -8151.7646
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H55").Value = 26896.053
$ws.Range("I55").Value = 50676.2
$ws.Range("J55").Value = 473.66666
$ws.Range("K55").Value = 50676.2
$ws.Range("L55").Value = 473.66666
$ws.Range("M55").Value = -50503.2
$ws.Range("N55").Value = -819.66666
$ws.Range("H63").Value = 58000
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 58000
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 72494.25
$ws.Range("I54").Value = 200000
$ws.Range("J54").Value = 29992.334
$ws.Range("K54").Value = 200000
$ws.Range("L54").Value = 29992.334
$ws.Range("M54").Value = -199480
$ws.Range("N54").Value = -31032.334
$ws.Range("H107").Value = 3250
$ws.Range("I107").Value = 2978.5833
$ws.Range("K107").Value = 8935.749899999999
$ws.Range("M107").Value = -7015.749899999999
$ws.Range("H126").Value = 3148.75
$ws.Range("I126").Value = 1689.6666
$ws.Range("K126").Value = 5068.9998
$ws.Range("M126").Value = -2598.9998
